$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 557.5714
$ws.Range("I6").Value = 557.5714
$ws.Range("K6").Value = 1672.7142
$ws.Range("M6").Value = -1560.7142
$ws.Range("H17").Value = 600.1539
$ws.Range("J17").Value = 608.2105
$ws.Range("L17").Value = 1824.6315
$ws.Range("N17").Value = -2160.6315
$ws.Range("H28").Value = 1573.6666
$ws.Range("I28").Value = 622.8333
$ws.Range("J28").Value = 3475.3333
$ws.Range("K28").Value = 622.8333
$ws.Range("L28").Value = 3475.3333
$ws.Range("M28").Value = -137.8333
$ws.Range("N28").Value = -4445.3333
$ws.Range("H53").Value = 471.9375
$ws.Range("I53").Value = 149.71428
$ws.Range("J53").Value = 1087.091
$ws.Range("K53").Value = 149.71428
$ws.Range("L53").Value = 1087.091
$ws.Range("M53").Value = 487.28572
$ws.Range("N53").Value = -2361.091
$ws.Range("H76").Value = 3397.25
$ws.Range("I76").Value = 2917.5
$ws.Range("K76").Value = 2917.5
$ws.Range("M76").Value = -2602.5
$ws.Range("H79").Value = 3397.25
$ws.Range("I79").Value = 2917.5
$ws.Range("K79").Value = 2917.5
$ws.Range("M79").Value = -1825.5
$ws.Range("H101").Value = 473.5
$ws.Range("J101").Value = 531.3333
$ws.Range("L101").Value = 1593.9999
$ws.Range("N101").Value = -4837.9999
$ws.Range("H112").Value = 5042.364
$ws.Range("J112").Value = 4718.4443
$ws.Range("L112").Value = 14155.3329
$ws.Range("N112").Value = -16371.3329
$ws.Range("H129").Value = 1088.6552
$ws.Range("H138").Value = 4742.846
$ws.Range("I138").Value = 4742.846
$ws.Range("K138").Value = 14228.538
$ws.Range("M138").Value = -9088.537999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 13417.083
$ws.Range("I61").Value = 13817.637
$ws.Range("J61").Value = 9011
$ws.Range("K61").Value = 13817.637
$ws.Range("L61").Value = 9011
$ws.Range("M61").Value = -13605.637
$ws.Range("N61").Value = -9435
$ws.Range("H110").Value = 2174.7058
$ws.Range("I110").Value = 1019
$ws.Range("K110").Value = 1019
$ws.Range("M110").Value = 1026
$ws.Range("H132").Value = 27131.3
$ws.Range("I132").Value = 29562.889
$ws.Range("K132").Value = 88688.667
$ws.Range("M132").Value = -86158.667
$ws.Range("H136").Value = 13417.083
$ws.Range("I136").Value = 13817.637
$ws.Range("J136").Value = 9011
$ws.Range("K136").Value = 41452.911
$ws.Range("L136").Value = 27033
$ws.Range("M136").Value = -38902.911
$ws.Range("N136").Value = -32133

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3259.6
$ws.Range("I105").Value = 3026.3
$ws.Range("K105").Value = 3026.3
$ws.Range("M105").Value = -1279.3
$ws.Range("H112").Value = 130146
$ws.Range("J112").Value = 130146
$ws.Range("L112").Value = 130146
$ws.Range("N112").Value = -133100
$ws.Range("H134").Value = 2706.55
$ws.Range("I134").Value = 2015.2354
$ws.Range("K134").Value = 6045.706200000001
$ws.Range("M134").Value = -3510.706200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 851.1429000000001
$ws.Range("I22").Value = 766.7
$ws.Range("J22").Value = 1062.25
$ws.Range("K22").Value = 766.7
$ws.Range("L22").Value = 1062.25
$ws.Range("M22").Value = -416.7
$ws.Range("N22").Value = -1762.25
$ws.Range("H58").Value = 54886.473
$ws.Range("I58").Value = 68521.8
$ws.Range("J58").Value = 3754
$ws.Range("K58").Value = 68521.8
$ws.Range("L58").Value = 3754
$ws.Range("M58").Value = -68318.8
$ws.Range("N58").Value = -4160
$ws.Range("H80").Value = 21999
$ws.Range("J80").Value = 21999
$ws.Range("L80").Value = 21999
$ws.Range("N80").Value = -24245
$ws.Range("H83").Value = 21999
$ws.Range("J83").Value = 21999
$ws.Range("L83").Value = 65997
$ws.Range("N83").Value = -77229
$ws.Range("H105").Value = 1101.7142
$ws.Range("I105").Value = 1101.7142
$ws.Range("K105").Value = 1101.7142
$ws.Range("M105").Value = 645.2858000000001
$ws.Range("H132").Value = 2152.3462
$ws.Range("I132").Value = 1652.3077
$ws.Range("K132").Value = 4956.9231
$ws.Range("M132").Value = -2426.9231
$ws.Range("H136").Value = 54886.473
$ws.Range("I136").Value = 68521.8
$ws.Range("J136").Value = 3754
$ws.Range("K136").Value = 205565.4
$ws.Range("L136").Value = 11262
$ws.Range("M136").Value = -203015.4
$ws.Range("N136").Value = -16362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 480.2857
$ws.Range("I86").Value = 423
$ws.Range("J86").Value = 623.5
$ws.Range("K86").Value = 1269
$ws.Range("L86").Value = 1870.5
$ws.Range("M86").Value = -83
$ws.Range("N86").Value = -4242.5
$ws.Range("H89").Value = 480.2857
$ws.Range("I89").Value = 423
$ws.Range("J89").Value = 623.5
$ws.Range("K89").Value = 3807
$ws.Range("L89").Value = 5611.5
$ws.Range("M89").Value = 2121
$ws.Range("N89").Value = -17467.5
$ws.Range("H105").Value = 9151.625
$ws.Range("J105").Value = 9151.625
$ws.Range("L105").Value = 27454.875
$ws.Range("N105").Value = -32696.875
$ws.Range("H114").Value = 1280.1875
$ws.Range("I114").Value = 1620.75
$ws.Range("K114").Value = 4862.25
$ws.Range("M114").Value = -1608.25
$ws.Range("H119").Value = 4012.4443
$ws.Range("I119").Value = 2014.125
$ws.Range("K119").Value = 6042.375
$ws.Range("M119").Value = -1204.375
$ws.Range("H120").Value = 15829.167
$ws.Range("I120").Value = 12995.2
$ws.Range("J120").Value = 29999
$ws.Range("K120").Value = 38985.60000000001
$ws.Range("L120").Value = 89997
$ws.Range("M120").Value = -34147.60000000001
$ws.Range("N120").Value = -99673
$ws.Range("H140").Value = 3572.8333
$ws.Range("I140").Value = 3572.8333
$ws.Range("K140").Value = 10718.4999
$ws.Range("M140").Value = -5538.499899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 31363.545
$ws.Range("J15").Value = 28499.9
$ws.Range("L15").Value = 28499.9
$ws.Range("N15").Value = -29075.9
$ws.Range("H70").Value = 6331.7
$ws.Range("I70").Value = 6625
$ws.Range("J70").Value = 6136.1665
$ws.Range("K70").Value = 6625
$ws.Range("L70").Value = 6136.1665
$ws.Range("M70").Value = -6355
$ws.Range("N70").Value = -6676.1665
$ws.Range("H73").Value = 6331.7
$ws.Range("I73").Value = 6625
$ws.Range("J73").Value = 6136.1665
$ws.Range("K73").Value = 6625
$ws.Range("L73").Value = 6136.1665
$ws.Range("M73").Value = -5689
$ws.Range("N73").Value = -8008.1665
$ws.Range("H81").Value = 31363.545
$ws.Range("J81").Value = 28499.9
$ws.Range("L81").Value = 28499.9
$ws.Range("N81").Value = -30495.9
$ws.Range("H84").Value = 31363.545
$ws.Range("J84").Value = 28499.9
$ws.Range("L84").Value = 85499.70000000001
$ws.Range("N84").Value = -95483.70000000001
$ws.Range("H132").Value = 46291.957
$ws.Range("I132").Value = 52905.8
$ws.Range("K132").Value = 158717.4
$ws.Range("M132").Value = -156187.4
$ws.Range("H139").Value = 35000
$ws.Range("J139").Value = 35000
$ws.Range("L139").Value = 35000
$ws.Range("N139").Value = -45280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 38598.75
$ws.Range("J42").Value = 38999.5
$ws.Range("L42").Value = 38999.5
$ws.Range("N42").Value = -40125.5
$ws.Range("H43").Value = 15595.8
$ws.Range("J43").Value = 18994.75
$ws.Range("L43").Value = 18994.75
$ws.Range("N43").Value = -19380.75
$ws.Range("H49").Value = 38598.75
$ws.Range("J49").Value = 38999.5
$ws.Range("L49").Value = 38999.5
$ws.Range("N49").Value = -39293.5
$ws.Range("H74").Value = 68833
$ws.Range("I74").Value = 68250
$ws.Range("K74").Value = 68250
$ws.Range("M74").Value = -67252
$ws.Range("H77").Value = 68833
$ws.Range("I77").Value = 68250
$ws.Range("K77").Value = 204750
$ws.Range("M77").Value = -199758
$ws.Range("H100").Value = 4548.1
$ws.Range("I100").Value = 3299.6667
$ws.Range("J100").Value = 5083.143
$ws.Range("K100").Value = 3299.6667
$ws.Range("L100").Value = 5083.143
$ws.Range("M100").Value = -2758.6667
$ws.Range("N100").Value = -6165.143
$ws.Range("H122").Value = 4299.033
$ws.Range("I122").Value = 3481.8333
$ws.Range("K122").Value = 10445.4999
$ws.Range("M122").Value = -7995.499899999999
$ws.Range("H132").Value = 57368.957
$ws.Range("I132").Value = 80243.69
$ws.Range("J132").Value = 5083.857
$ws.Range("K132").Value = 240731.07
$ws.Range("L132").Value = 15251.571
$ws.Range("M132").Value = -238201.07
$ws.Range("N132").Value = -20311.571
$ws.Range("H133").Value = 94999.5
$ws.Range("J133").Value = 94999.5
$ws.Range("L133").Value = 94999.5
$ws.Range("N133").Value = -100059.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1433.1578
$ws.Range("I122").Value = 1320.6875
$ws.Range("J122").Value = 2033
$ws.Range("K122").Value = 3962.0625
$ws.Range("L122").Value = 6099
$ws.Range("M122").Value = -1512.0625
$ws.Range("N122").Value = -10999
$ws.Range("H126").Value = 34053.332
$ws.Range("I126").Value = 40553.703
$ws.Range("J126").Value = 4801.6665
$ws.Range("K126").Value = 121661.109
$ws.Range("L126").Value = 14404.9995
$ws.Range("M126").Value = -119191.109
$ws.Range("N126").Value = -19344.9995
$ws.Range("H132").Value = 29100.703
$ws.Range("I132").Value = 30657.885
$ws.Range("K132").Value = 91973.655
$ws.Range("M132").Value = -89443.655
